$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test rows 7-9 ("site" / "site privado" test cases) ---
# Values are entered in the same order the original author typed them so the
# shared-string table comes out in the same sequence as the saved workbook.
$ws.Range("A8").Value = "site"
$ws.Range("B7").Value = "Tentar mudar de pagina"
$ws.Range("C7").Value = "ao clicar em algum link, mudar para a pagina esperada"
$ws.Range("A7").Value = "site privado"
$ws.Range("B8").Value = "testar botões do menu"
$ws.Range("C8").Value = "ao clicar nos botões, ir para a parte correspondente da pagina"
$ws.Range("B9").Value = "testar os simuladores financeiros"
$ws.Range("C9").Value = 'ao clicar em "simular", em cada uma das opções, deve ser exibido um alert com o resultado'
$ws.Range("A9").Value = "site"
$ws.Range("D7").Value = "aguardando "
$ws.Range("D8").Value = "aguardando "
$ws.Range("D9").Value = "aguardando "

# D7:D9 are brand new cells - copy the existing row formatting (centered,
# wrapped text) from column A of the same row onto them.
$ws.Range("A7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("A8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 7-9 use the taller (30pt) row height, matching rows 4-6 above them
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30

# --- D10: new empty placeholder cell, underlined ---
$ws.Range("D10").Font.Underline = $true

# --- Cursor / selection moved from C12 to D10 ---
$ws.Range("D10").Select()
